$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that sometimes look like plain decimals
# (e.g. "234.01"); force text format so Excel does not coerce them to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.849.48"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.93"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.01"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.66"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.401.93"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.34"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.37"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.096.49"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.814.41"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.16"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.45"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.83"
$ws.Range("E26").Value = "  +8.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.35"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.54"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.40"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  +9.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.62"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0974"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  +4.94%  "
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.452.07"
$ws.Range("E47").Value = "  -5.25%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.285.53"
$ws.Range("E51").Value = "  +0.23%  "
